$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.485.20"
$ws.Range("D3").Value = "1.838.10"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.46"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5369"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2913"
$ws.Range("E8").Value = "  -9.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06958"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.18"
$ws.Range("E10").Value = "  -9.27%  "
$ws.Range("D11").Value = "1.847.39"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7236"
$ws.Range("E12").Value = "  -7.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07191"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.96"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.973"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.75"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007886"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").Value = "26.501.69"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "2.082.93"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.579"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.979"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.170"
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.87"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.152"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.706"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.93"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.77"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.244"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08884"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04836"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.898"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7225"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.129"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.090"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.291"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01706"
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4651"
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9019"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.81"
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.861"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.380"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.993"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1240"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.74"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4035"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8902"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05744"
$ws.Range("E51").Value = "  -2.29%  "
